# Update "want to go" counts (column F) for the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 7592
    5  = 323
    6  = 31
    7  = 20
    8  = 23
    9  = 5723
    10 = 145
    11 = 10
    12 = 18
    13 = 1754
    14 = 60
    15 = 1241
    16 = 277
    17 = 5501
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
